$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.176
$ws.Range("C3").Value = -13.09
$ws.Range("C14").Value = -12.509
$ws.Range("C21").Value = -12.499
$ws.Range("C23").Value = -12.594
$ws.Range("C25").Value = -12.699
$ws.Range("D25").Value = -8.424000000000001
$ws.Range("C26").Value = -13.117
$ws.Range("D27").Value = -8.597000000000001
$ws.Range("C29").Value = -12.2
$ws.Range("D31").Value = -8.394
$ws.Range("D39").Value = -7.846000000000001
$ws.Range("D48").Value = -7.475
$ws.Range("D51").Value = -8.373999999999999
$ws.Range("D52").Value = -7.469000000000001
$ws.Range("C53").Value = -11.04
$ws.Range("D55").Value = -8.065
$ws.Range("D56").Value = -8.272
$ws.Range("C57").Value = -13.565
$ws.Range("D57").Value = -8.559000000000001
$ws.Range("C59").Value = -13.077
$ws.Range("C69").Value = -10.676
$ws.Range("D73").Value = -8.004000000000001
$ws.Range("C79").Value = -12.491
$ws.Range("C83").Value = -13.169
$ws.Range("D89").Value = -6.858
$ws.Range("D90").Value = -7.601999999999999
$ws.Range("C91").Value = -10.571
$ws.Range("D92").Value = -6.593000000000001
$ws.Range("C93").Value = -11.886
